# Weekly update for "Poroto verde" - Macroferia Regional de Talca.
# A new weekly record is inserted (previous row 155 data is preserved as a
# new row), and row 155 is updated with this week's new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 156 downward to make room for the duplicated
# (previous week's) record; the old row 156 becomes row 157 unchanged.
$ws.Rows.Item(156).Insert()

# Re-create the record that used to live in row 155 (before this week's
# update) now at row 156.
$ws.Cells.Item(156, 1).Value = 5
$ws.Cells.Item(156, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(156, 3).Value = "Maule"
$ws.Cells.Item(156, 4).Value = 44508
$ws.Cells.Item(156, 5).Value = 7
$ws.Cells.Item(156, 6).Value = 100112031
$ws.Cells.Item(156, 7).Value = "Poroto verde"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 150
$ws.Cells.Item(156, 11).Value = 40000
$ws.Cells.Item(156, 12).Value = 40000
$ws.Cells.Item(156, 13).Value = 40000
$ws.Cells.Item(156, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(156, 15).Value = "Región del Maule"
$ws.Cells.Item(156, 16).Value = 1600
$ws.Cells.Item(156, 17).Value = 25
$ws.Cells.Item(156, 18).Value = "Hortaliza"

# Update row 155 with this week's new price/volume data.
$ws.Cells.Item(155, 4).Value = 44656
$ws.Cells.Item(155, 10).Value = 100
$ws.Cells.Item(155, 11).Value = 25000
$ws.Cells.Item(155, 12).Value = 25000
$ws.Cells.Item(155, 13).Value = 25000
$ws.Cells.Item(155, 16).Value = 1000
